$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the selected/active cell on the sheet view (scroll back to top-left,
# select M5 instead of Z25)
$ws.Range("M5").Select()

# Swap the clue-label text values that were corrected (room center / room
# label swaps) by writing the new cell values directly.
$ws.Range("E4").Value = "A*"
$ws.Range("O4").Value = "T*"
$ws.Range("D5").Value = "A#"
$ws.Range("M5").Value = "T#"
$ws.Range("X5").Value = "C*"
$ws.Range("W6").Value = "C#"
$ws.Range("X10").Value = "J*"
$ws.Range("W12").Value = "J#"
$ws.Range("E14").Value = "N*"
$ws.Range("B16").Value = "N#"
$ws.Range("V16").Value = "U#"
$ws.Range("Y16").Value = "U"
$ws.Range("V17").Value = "U"
$ws.Range("X17").Value = "U*"
$ws.Range("X24").Value = "Q*"
$ws.Range("D25").Value = "S*"
$ws.Range("N25").Value = "D#"
$ws.Range("V26").Value = "Q#"
$ws.Range("C27").Value = "S#"
$ws.Range("N27").Value = "D*"

$wb.Save()
